$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9626310369308592
$ws.Range("D2").Value = 0.9886668653317213
$ws.Range("E2").Value = 0.9717025505326101
$ws.Range("F2").Value = 0.9373878623774127
$ws.Range("I2").Value = 1.030769576117298
$ws.Range("J2").Value = 0.9859899374366492
$ws.Range("K2").Value = 1.00028637980606
$ws.Range("L2").Value = 0.9835707165988693
$ws.Range("M2").Value = 0.9497943589845799
$ws.Range("C3").Value = 0.9741197046967054
$ws.Range("D3").Value = 0.9975400573709011
$ws.Range("E3").Value = 0.9816447563085319
$ws.Range("F3").Value = 0.9521367267359822
$ws.Range("I3").Value = 1.034760932073159
$ws.Range("J3").Value = 0.9953416490766156
$ws.Range("K3").Value = 1.008204510688724
$ws.Range("L3").Value = 0.9925180111125878
$ws.Range("M3").Value = 0.9634210127105038
$ws.Range("C4").Value = 0.9811833483245908
$ws.Range("D4").Value = 1.002998958544201
$ws.Range("E4").Value = 0.9877593571912343
$ws.Range("F4").Value = 0.9611880062170053
$ws.Range("I4").Value = 1.037195021052741
$ws.Range("J4").Value = 1.001080948628347
$ws.Range("K4").Value = 1.01306132303485
$ws.Range("L4").Value = 0.9980077960741711
$ws.Range("M4").Value = 0.9717785839467528
$ws.Range("C5").Value = 0.9840725915156856
$ws.Range("D5").Value = 1.005232356364963
$ws.Range("E5").Value = 0.9902606757668951
$ws.Range("F5").Value = 0.964887236458749
$ws.Range("I5").Value = 1.03818575241667
$ws.Range("J5").Value = 1.003425946718123
$ws.Range("K5").Value = 1.015045012812038
$ws.Range("L5").Value = 1.000250493367595
$ws.Range("M5").Value = 0.9751929860047728
$ws.Range("C6").Value = 0.9845532168572377
$ws.Range("D6").Value = 1.005603906852584
$ws.Range("E6").Value = 0.9906767809865954
$ws.Range("F6").Value = 0.9655024521337245
$ws.Range("I6").Value = 1.038350270176514
$ws.Range("J6").Value = 1.00381588585945
$ws.Range("K6").Value = 1.015374825721242
$ws.Range("L6").Value = 1.000623399680374
$ws.Range("M6").Value = 0.9757607521484258
$ws.Range("C7").Value = 0.9812222600595656
$ws.Range("D7").Value = 1.003029035653499
$ws.Range("E7").Value = 0.9877930435948024
$ws.Range("F7").Value = 0.9612378374439015
$ws.Range("I7").Value = 1.03720838355773
$ws.Range("J7").Value = 1.001112540757652
$ws.Range("K7").Value = 1.013088050615052
$ws.Range("L7").Value = 0.9980380114385742
$ws.Range("M7").Value = 0.9718245836818195
$ws.Range("C8").Value = 0.9665958218229636
$ws.Range("D8").Value = 0.9917281491470212
$ws.Range("E8").Value = 0.975133154110868
$ws.Range("F8").Value = 0.942482067176491
$ws.Range("I8").Value = 1.032151212267571
$ws.Range("J8").Value = 0.989219467582792
$ws.Range("K8").Value = 1.003021346131679
$ws.Range("L8").Value = 0.9866608417472876
$ws.Range("M8").Value = 0.9545019562034169
$ws.Range("C9").Value = 0.9375207092744419
$ws.Range("D9").Value = 0.9693075607549898
$ws.Range("E9").Value = 0.9499922893967362
$ws.Range("F9").Value = 0.9049895117027649
$ws.Range("I9").Value = 1.021933945508186
$ws.Range("J9").Value = 0.965490555173467
$ws.Range("K9").Value = 0.982919036116091
$ws.Range("L9").Value = 0.963952394368968
$ws.Range("M9").Value = 0.919838312084708
$ws.Range("C10").Value = 0.9150576992282725
$ws.Range("D10").Value = 0.9520442931450224
$ws.Range("E10").Value = 0.9306048912318564
$ws.Range("F10").Value = 0.8757457548717932
$ws.Range("I10").Value = 1.013934927751297
$ws.Range("J10").Value = 0.9470997962636374
$ws.Range("K10").Value = 0.9673369685553551
$ws.Range("L10").Value = 0.9463513184561065
$ws.Range("M10").Value = 0.8927903691995348
$ws.Range("C11").Value = 0.904314053425991
$ws.Range("D11").Value = 0.9438102772121674
$ws.Range("E11").Value = 0.9213468645887903
$ws.Range("F11").Value = 0.861645833607909
$ws.Range("I11").Value = 1.010086775679746
$ws.Range("J11").Value = 0.938290566387191
$ws.Range("K11").Value = 0.9598757867404154
$ws.Range("L11").Value = 0.9379218759397218
$ws.Range("M11").Value = 0.8797518762683875
$ws.Range("B12").Value = 1.05
$ws.Range("C12").Value = 0.9001349369169596
$ws.Range("D12").Value = 0.940611772245459
$ws.Range("E12").Value = 0.9177485320943363
$ws.Range("F12").Value = 0.8561381416992523
$ws.Range("I12").Value = 1.008586859281919
$ws.Range("J12").Value = 0.9348619635728417
$ws.Range("K12").Value = 0.9569726252708672
$ws.Range("L12").Value = 0.9346415554967482
$ws.Range("M12").Value = 0.8746598972223139
$ws.Range("C13").Value = 0.9010406193240276
$ws.Range("D13").Value = 0.9413047234826143
$ws.Range("E13").Value = 0.9185282052469104
$ws.Range("F13").Value = 0.857332905631905
$ws.Range("I13").Value = 1.008912046242702
$ws.Range("J13").Value = 0.9356050842515505
$ws.Range("K13").Value = 0.9576018190672825
$ws.Range("L13").Value = 0.9353525113948766
$ws.Range("M13").Value = 0.8757644185486632
$ws.Range("C14").Value = 0.9039727046874908
$ws.Range("D14").Value = 0.9435489293815931
$ws.Range("E14").Value = 0.921052891037881
$ws.Range("F14").Value = 0.8611964729078219
$ws.Range("I14").Value = 1.00996432205915
$ws.Range("J14").Value = 0.9380105581811979
$ws.Range("K14").Value = 0.9596386723937698
$ws.Range("L14").Value = 0.9376539666309047
$ws.Range("M14").Value = 0.8793364056960195
$ws.Range("C15").Value = 0.9057529973430571
$ws.Range("D15").Value = 0.9449121634526536
$ws.Range("E15").Value = 0.9225862212197896
$ws.Range("F15").Value = 0.8635391315464301
$ws.Range("I15").Value = 1.010602854577959
$ws.Range("J15").Value = 0.9394708535209874
$ws.Range("K15").Value = 0.9608753017309138
$ws.Range("L15").Value = 0.9390511843418023
$ws.Range("M15").Value = 0.8815024332814811
$ws.Range("C16").Value = 0.9157465301262258
$ws.Range("D16").Value = 0.9525727523430371
$ws.Range("E16").Value = 0.9311988196248185
$ws.Range("F16").Value = 0.8766470444721205
$ws.Range("I16").Value = 1.014181225537159
$ws.Range("J16").Value = 0.9476643381514307
$ws.Range("K16").Value = 0.9678152024677559
$ws.Range("L16").Value = 0.9468915698607634
$ws.Range("M16").Value = 0.8936239158667587
$ws.Range("C17").Value = 0.9217214422396119
$ws.Range("D17").Value = 0.9571591705515647
$ws.Range("E17").Value = 0.9363522062296105
$ws.Range("F17").Value = 0.8844519850694066
$ws.Range("I17").Value = 1.016315193185388
$ws.Range("J17").Value = 0.9525597252145025
$ws.Range("K17").Value = 0.9719625000241047
$ws.Range("L17").Value = 0.9515765084632948
$ws.Range("M17").Value = 0.9008425671779018
$ws.Range("C18").Value = 0.9251113108309676
$ws.Range("D18").Value = 0.959763262995983
$ws.Range("E18").Value = 0.9392772558053692
$ws.Range("F18").Value = 0.8888703122593442
$ws.Range("I18").Value = 1.01752384823012
$ws.Range("J18").Value = 0.9553359184725286
$ws.Range("K18").Value = 0.974314664245293
$ws.Range("L18").Value = 0.9542334758284896
$ws.Range("M18").Value = 0.9049291921259247
$ws.Range("C19").Value = 0.9262519110647932
$ws.Range("D19").Value = 0.9606397832310732
$ws.Range("E19").Value = 0.9402616559604619
$ws.Range("F19").Value = 0.8903554355946184
$ws.Range("I19").Value = 1.017930172293869
$ws.Range("J19").Value = 0.9562698271110031
$ws.Range("K19").Value = 0.9751059575394643
$ws.Range("L19").Value = 0.955127291972524
$ws.Range("M19").Value = 0.9063028380127247
$ws.Range("C20").Value = 0.9210904556031323
$ws.Range("D20").Value = 0.9566746035879298
$ws.Range("E20").Value = 0.9358078402760385
$ws.Range("F20").Value = 0.8836287961697028
$ws.Range("I20").Value = 1.016090045949251
$ws.Range("J20").Value = 0.9520428686403599
$ws.Range("K20").Value = 0.9715246022460836
$ws.Range("L20").Value = 0.9510818571814417
$ws.Range("M20").Value = 0.9000811903023016
$ws.Range("C21").Value = 0.9031148339939379
$ws.Range("D21").Value = 0.9428921894768314
$ws.Range("E21").Value = 0.9203141317649017
$ws.Range("F21").Value = 0.8600667564151723
$ws.Range("I21").Value = 1.009656526767692
$ws.Range("J21").Value = 0.9373068163578632
$ws.Range("K21").Value = 0.9590427490767796
$ws.Range("L21").Value = 0.936980641271825
$ws.Range("M21").Value = 0.8782919111358197
$ws.Range("C22").Value = 0.8906935541877636
$ws.Range("D22").Value = 0.9333953991352029
$ws.Range("E22").Value = 0.9096256800999449
$ws.Range("F22").Value = 0.843642924808003
$ws.Range("I22").Value = 1.005193023972452
$ws.Range("J22").Value = 0.9271125255348823
$ws.Range("K22").Value = 0.9504127854520259
$ws.Range("L22").Value = 0.9272285095471251
$ws.Range("M22").Value = 0.863110933009245
$ws.Range("C23").Value = 0.8974002148350931
$ws.Range("D23").Value = 0.9385201349741438
$ws.Range("E23").Value = 0.9153947901448034
$ws.Range("F23").Value = 0.8525265378640369
$ws.Range("I23").Value = 1.007604535261551
$ws.Range("J23").Value = 0.9326178164204675
$ws.Range("K23").Value = 0.9550726745602265
$ws.Range("L23").Value = 0.9324946354860215
$ws.Range("M23").Value = 0.8713213164869571
$ws.Range("C24").Value = 0.9213758629571591
$ws.Range("D24").Value = 0.9568937763752323
$ws.Range("E24").Value = 0.9360540635005281
$ws.Range("F24").Value = 0.8840011705534472
$ws.Range("I24").Value = 1.016191890798392
$ws.Range("J24").Value = 0.9522766565539061
$ws.Range("K24").Value = 0.9717226743570925
$ws.Range("L24").Value = 0.9513056006931294
$ws.Range("M24").Value = 0.9004256030938486
$ws.Range("C25").Value = 0.9455038186231672
$ws.Range("D25").Value = 0.9754561827162872
$ws.Range("E25").Value = 0.9568907735382788
$ws.Range("F25").Value = 0.9153181483625437
$ws.Range("I25").Value = 1.024757229012188
$ws.Range("J25").Value = 0.972015468293208
$ws.Range("K25").Value = 0.9884478466182575
$ws.Range("L25").Value = 0.9701973222565364
$ws.Range("M25").Value = 0.92939058333222
